# Apply the cibmtr-reporting IG metadata update to the ValueSet-presence-valueset workbook.

$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

# Version bump
$meta.Range("B3").Value = "0.1.7"

# Status: active -> draft
$meta.Range("B6").Value = "draft"

# Date update
$meta.Range("B8").Value = "2024-08-27T12:23:18-05:00"

# Contact row 10 now carries the real publisher contact detail (with URL).
$meta.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# Contact row 11 (previously a duplicate of row 10) becomes the second named contact.
$meta.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# Insert a new row 12 for "Jurisdiction" (empty value), shifting Description,
# Purpose, Copyright and Immutable down by one row.
$meta.Rows("12").Insert()
$meta.Range("A12").Value = "Jurisdiction"
$meta.Range("B12").Value = ""
$meta.Range("A12:B12").Style = $meta.Range("A11:B11").Style
